$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 537 (shifts old rows 537..563 down to 538..564,
# and carries the existing formatting, e.g. the date style on column D).
$ws.Rows.Item(537).Insert()

# Populate the newly inserted row 537 with the new weekly record.
$ws.Cells.Item(537, 1).Value = 9
$ws.Cells.Item(537, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(537, 3).Value = "Metropolitana"
$ws.Cells.Item(537, 4).Value = (Get-Date -Year 2022 -Month 7 -Day 11 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(537, 5).Value = 13
$ws.Cells.Item(537, 6).Value = 100114014
$ws.Cells.Item(537, 7).Value = "Betarraga"
$ws.Cells.Item(537, 8).Value = "Sin especificar"
$ws.Cells.Item(537, 9).Value = "Primera"
$ws.Cells.Item(537, 10).Value = 4300
$ws.Cells.Item(537, 11).Value = 150
$ws.Cells.Item(537, 12).Value = 160
$ws.Cells.Item(537, 13).Value = 155
$ws.Cells.Item(537, 14).Value = "$/unidad"
$ws.Cells.Item(537, 15).Value = "Región Metropolitana"
$ws.Cells.Item(537, 16).Value = 155
$ws.Cells.Item(537, 17).Value = 1
$ws.Cells.Item(537, 18).Value = "Hortaliza"

Write-Output "done"
